$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 481, shifting existing rows 481:523 down to 482:524
$ws.Rows.Item(481).Insert()

# Populate the newly inserted row 481 with this week's data
$ws.Cells.Item(481, 1).Value = 4
$ws.Cells.Item(481, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(481, 3).Value = "Los Lagos"
$ws.Cells.Item(481, 4).Value = 45166
$ws.Cells.Item(481, 5).Value = 10
$ws.Cells.Item(481, 6).Value = 100112017
$ws.Cells.Item(481, 7).Value = "Apio"
$ws.Cells.Item(481, 8).Value = "Americana (o)"
$ws.Cells.Item(481, 9).Value = "Primera"
$ws.Cells.Item(481, 10).Value = 15
$ws.Cells.Item(481, 11).Value = 11000
$ws.Cells.Item(481, 12).Value = 11000
$ws.Cells.Item(481, 13).Value = 11000
$ws.Cells.Item(481, 14).Value = "`$/docena de matas"
$ws.Cells.Item(481, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(481, 16).Value = 1833
$ws.Cells.Item(481, 17).Value = 6
$ws.Cells.Item(481, 18).Value = "Hortaliza"
